$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds strikeouts (header "K"); regenerate save_data to use K
# instead of Strike# for rows 2-29 (data rows), per regen of std/mean and
# s_vals calculation.
$newK = @{
    2  = 6
    3  = 0
    4  = 3
    5  = 5
    6  = 6
    7  = 6
    8  = 3
    9  = 2
    10 = 5
    11 = 1
    12 = 0
    13 = 2
    14 = 7
    15 = 3
    16 = 3
    17 = 5
    18 = 8
    19 = 2
    20 = 4
    21 = 4
    22 = 3
    23 = 4
    24 = 1
    25 = 2
    26 = 3
    27 = 5
    28 = 4
    29 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
